$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.668.19"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.579.95"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.91%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.65"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.92"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.43%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.528"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.576.50"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.140"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.17"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.81"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.043.81"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000178"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.468.42"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.555.97"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.43"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -6.23%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "351.52"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.24"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.62"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.91"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "69.20"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.88"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -9.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.711.49"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.67%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0993"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "531.65"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.25%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.54%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.01%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.82%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.47"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "156.57"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.80"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.360"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.33"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.23%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.14"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.46"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.95%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "149.73"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.27%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₆0286"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.568"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.73"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.73"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0764"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.91%  "
